$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new question row (row 14): a 13th question ("Q13") under the
# "Énergie" axis, mirroring the other rows in the table.
$ws.Range("A14").Value = "Énergie"
$ws.Range("B14").Value = "Q13"
$ws.Range("C14").Value = "test insertion question"
$ws.Range("D14").Value = 1

# Match the look of the existing data rows (wrap text, vertically centered)
# by copying the formatting already used on row 13 instead of building a
# one-off style.
$ws.Range("A13:C13").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection on the newly entered cell, like Excel does right
# after typing a value in.
$ws.Range("D14").Select()
